$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 17571.428
$ws.Range("I32").Value = 28999
$ws.Range("J32").Value = 13000.4
$ws.Range("K32").Value = 28999
$ws.Range("L32").Value = 13000.4
$ws.Range("M32").Value = -28673
$ws.Range("N32").Value = -13652.4

$ws.Range("H40").Value = 3459.4614
$ws.Range("I40").Value = 3097.5557
$ws.Range("K40").Value = 3097.5557
$ws.Range("M40").Value = -2922.5557

$ws.Range("H137").Value = 6508.1055
$ws.Range("I137").Value = 4199.6665
$ws.Range("J137").Value = 12174.272
$ws.Range("K137").Value = 12598.9995
$ws.Range("L137").Value = 36522.81600000001
$ws.Range("M137").Value = -10048.9995
$ws.Range("N137").Value = -41622.81600000001

$ws.Range("H138").Value = 2951.9207
$ws.Range("I138").Value = 2717.375
$ws.Range("J138").Value = 4828.2856
$ws.Range("K138").Value = 8152.125
$ws.Range("L138").Value = 14484.8568
$ws.Range("M138").Value = -3012.125
$ws.Range("N138").Value = -24764.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2047231.6
$ws.Range("I2").Value = 2302678
$ws.Range("J2").Value = 3659.8
$ws.Range("K2").Value = 2302678
$ws.Range("L2").Value = 3659.8
$ws.Range("M2").Value = -2302565
$ws.Range("N2").Value = -3885.8

$ws.Range("H97").Value = 1326578.8
$ws.Range("I97").Value = 1767448
$ws.Range("K97").Value = 1767448
$ws.Range("M97").Value = -1766952

$ws.Range("H116").Value = 2047231.6
$ws.Range("I116").Value = 2302678
$ws.Range("J116").Value = 3659.8
$ws.Range("K116").Value = 2302678
$ws.Range("L116").Value = 3659.8
$ws.Range("M116").Value = -2300384
$ws.Range("N116").Value = -8247.799999999999

$ws.Range("H122").Value = 1396.5
$ws.Range("I122").Value = 1396.5
$ws.Range("K122").Value = 4189.5
$ws.Range("M122").Value = -1739.5

$ws.Range("H132").Value = 3674.194
$ws.Range("I132").Value = 2714.5518
$ws.Range("J132").Value = 9858.556
$ws.Range("K132").Value = 8143.655400000001
$ws.Range("L132").Value = 29575.668
$ws.Range("M132").Value = -5613.655400000001
$ws.Range("N132").Value = -34635.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2047231.6
$ws.Range("I3").Value = 2302678
$ws.Range("J3").Value = 3659.8
$ws.Range("K3").Value = 2302678
$ws.Range("L3").Value = 3659.8
$ws.Range("M3").Value = -2302564
$ws.Range("N3").Value = -3887.8

$ws.Range("H134").Value = 4365.0225
$ws.Range("I134").Value = 2964.4
$ws.Range("J134").Value = 7166.2666
$ws.Range("K134").Value = 8893.200000000001
$ws.Range("L134").Value = 21498.7998
$ws.Range("M134").Value = -6358.200000000001
$ws.Range("N134").Value = -26568.7998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4324.58
$ws.Range("I31").Value = 2173.5715
$ws.Range("J31").Value = 5161.0835
$ws.Range("K31").Value = 2173.5715
$ws.Range("L31").Value = 5161.0835
$ws.Range("M31").Value = -1878.5715
$ws.Range("N31").Value = -5751.0835

$ws.Range("H34").Value = 4324.58
$ws.Range("I34").Value = 2173.5715
$ws.Range("J34").Value = 5161.0835
$ws.Range("K34").Value = 2173.5715
$ws.Range("L34").Value = 5161.0835
$ws.Range("M34").Value = -1971.5715
$ws.Range("N34").Value = -5565.0835

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 152.36363
$ws.Range("J12").Value = 157.7
$ws.Range("L12").Value = 473.1
$ws.Range("N12").Value = -819.0999999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 50001
$ws.Range("J93").Value = 50001
$ws.Range("L93").Value = 50001
$ws.Range("N93").Value = -53745

$ws.Range("H102").Value = 1858.9231
$ws.Range("I102").Value = 2154.0386
$ws.Range("J102").Value = 1563.8077
$ws.Range("K102").Value = 2154.0386
$ws.Range("L102").Value = 1563.8077
$ws.Range("M102").Value = -532.0385999999999
$ws.Range("N102").Value = -4807.8077

$ws.Range("H113").Value = 15879.083
$ws.Range("I113").Value = 21337.375
$ws.Range("J113").Value = 4962.5
$ws.Range("K113").Value = 21337.375
$ws.Range("L113").Value = 4962.5
$ws.Range("M113").Value = -19167.375
$ws.Range("N113").Value = -9302.5

$ws.Range("H126").Value = 3138.261
$ws.Range("I126").Value = 2283.2307
$ws.Range("J126").Value = 4249.8
$ws.Range("K126").Value = 6849.6921
$ws.Range("L126").Value = 12749.4
$ws.Range("M126").Value = -4379.6921
$ws.Range("N126").Value = -17689.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 64698.562
$ws.Range("I7").Value = 85098.164
$ws.Range("J7").Value = 3499.75
$ws.Range("K7").Value = 85098.164
$ws.Range("L7").Value = 3499.75
$ws.Range("M7").Value = -84986.164
$ws.Range("N7").Value = -3723.75

$ws.Range("H40").Value = 14711918
$ws.Range("I40").Value = 16134336
$ws.Range("K40").Value = 16134336
$ws.Range("M40").Value = -16134200

$ws.Range("H56").Value = 16316.667
$ws.Range("I56").Value = 11580
$ws.Range("J56").Value = 40000
$ws.Range("K56").Value = 11580
$ws.Range("L56").Value = 40000
$ws.Range("M56").Value = -10889
$ws.Range("N56").Value = -41382

$ws.Range("H61").Value = 1370.5385
$ws.Range("I61").Value = 856.5
$ws.Range("J61").Value = 3084
$ws.Range("K61").Value = 856.5
$ws.Range("L61").Value = 3084
$ws.Range("M61").Value = -654.5
$ws.Range("N61").Value = -3488

$ws.Range("H113").Value = 1370.5385
$ws.Range("I113").Value = 856.5
$ws.Range("J113").Value = 3084
$ws.Range("K113").Value = 856.5
$ws.Range("L113").Value = 3084
$ws.Range("M113").Value = 1313.5
$ws.Range("N113").Value = -7424

$ws.Range("H126").Value = 64698.562
$ws.Range("I126").Value = 85098.164
$ws.Range("J126").Value = 3499.75
$ws.Range("K126").Value = 255294.492
$ws.Range("L126").Value = 10499.25
$ws.Range("M126").Value = -252824.492
$ws.Range("N126").Value = -15439.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 49976.668
$ws.Range("J64").Value = 49976.668
$ws.Range("L64").Value = 49976.668
$ws.Range("N64").Value = -50472.668

$ws.Range("H67").Value = 49976.668
$ws.Range("J67").Value = 49976.668
$ws.Range("L67").Value = 49976.668
$ws.Range("N67").Value = -51692.668

$ws.Range("H122").Value = 1503.9714
$ws.Range("I122").Value = 1503.9714
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4511.914199999999
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -2061.914199999999

$ws.Range("H126").Value = 4635.12
$ws.Range("I126").Value = 4588.6665
$ws.Range("J126").Value = 4754.5713
$ws.Range("K126").Value = 13765.9995
$ws.Range("L126").Value = 14263.7139
$ws.Range("M126").Value = -11295.9995
$ws.Range("N126").Value = -19203.7139
